$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data
$ws.Range("D2").Value = '27.660.18'
$ws.Range("E2").Value = '  +3.17%  '
$ws.Range("D3").Value = '1.853.79'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.033'
$ws.Range("E4").Value = '  +2.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.71'
$ws.Range("E5").Value = '  +4.10%  '
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4382'
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3760'
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07422'
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8768'
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.51'
$ws.Range("E11").Value = '  +3.32%  '
$ws.Range("D12").Value = '1.868.86'
$ws.Range("E12").Value = '  -5.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.521'
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.708'
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07196'
$ws.Range("E15").Value = '  +4.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.01'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.034'
$ws.Range("E17").Value = '  +2.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009054'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '27.672.63'
$ws.Range("E21").Value = '  +3.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.277'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").Value = '2.077.16'
$ws.Range("E24").Value = '  -6.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.55'
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.946'
$ws.Range("E26").Value = '  +3.79%  '
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.297'
$ws.Range("E28").Value = '  +1.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.938'
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.46'
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09089'
$ws.Range("E31").Value = '  +1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.209'
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7698'
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.524'
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.883'
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.156'
$ws.Range("E37").Value = '  +1.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01983'
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05295'
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.827'
$ws.Range("E40").Value = '  +6.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5187'
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1676'
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.738'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.585'
$ws.Range("E44").Value = '  +3.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.01'
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.63'
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.721'
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4664'
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06398'
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.892'
$ws.Range("E50").Value = '  +4.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.63'
$ws.Range("E51").Value = '  +5.71%  '
